$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new column at E (pushes the old E -> F and old F -> G), mirroring
# a new "No AGDU" model-comparison column being added to the stats table.
# ---------------------------------------------------------------------------
$ws.Columns("E:E").Insert() | Out-Null

# Re-apply an explicit width to the freshly inserted column E so it lines up
# with the existing C:D columns (~18.1 chars wide) instead of the default.
$ws.Columns("E:E").ColumnWidth = 17.33

# ---------------------------------------------------------------------------
# Populate the new column E with the "No AGDU" comparison values.
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "No AGDU"

$ws.Range("E3").Value = 2.807
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 0.094
$ws.Range("E6").Value = 2.807

$ws.Range("E8").Value = 25.817
$ws.Range("E9").Value = 6
$ws.Range("E10").Value = 0

$ws.Range("E12").Value = 0.909
$ws.Range("E13").Value = 0.423

$ws.Range("E15").Value = 605.152
$ws.Range("E16").Value = 622.927
$ws.Range("E17").Value = 582.479

$ws.Range("E19").Value = 0.25
$ws.Range("E20").Value = 0.106

$ws.Range("E24").Value = 0.445
$ws.Range("E25").Value = 0.804
$ws.Range("E26").Value = 0.422

$ws.Range("E32").Value = 1.101
$ws.Range("E33").Value = 1.433
$ws.Range("E34").Value = 0.152

$ws.Range("E40").Value = 0.566
$ws.Range("E41").Value = 3.362
$ws.Range("E42").Value = 0.001

$ws.Range("E44").Value = 0.369
$ws.Range("E45").Value = 1.551
$ws.Range("E46").Value = 0.121

$ws.Range("E48").Value = 0.407
$ws.Range("E49").Value = 0.759
$ws.Range("E50").Value = 0.448

$ws.Range("E56").Value = 0.403
$ws.Range("E57").Value = 0.471
$ws.Range("E58").Value = 0.438

$ws.Range("E60").Value = 1.164
$ws.Range("E61").Value = 0
$ws.Range("E62").Value = 1.23
$ws.Range("E63").Value = 0

$ws.Range("E68").Value = 1.406
$ws.Range("E69").Value = 0

# ---------------------------------------------------------------------------
# Freeze the first column and move the active selection to E70, matching the
# new sheet view after the edit.
# ---------------------------------------------------------------------------
$ws.Range("B1").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E70").Select() | Out-Null
